$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("U4").Value = 2.38
$ws.Range("V4").Value = 1.53
$ws.Range("Y4").Value = 9.5
$ws.Range("AC4").Value = 6
$ws.Range("AG4").Value = 9.5
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 5.75
$ws.Range("X5").Value = 7
$ws.Range("Z5").Value = 11
$ws.Range("AB5").Value = 29
$ws.Range("AG5").Value = 15
$ws.Range("AH5").Value = 29
$ws.Range("AJ5").Value = 67
$ws.Range("AM5").Value = 351
$ws.Range("AN5").Value = 3.5
$ws.Range("AO5").Value = 8
$ws.Range("G7").Value = 2.18
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.77
$ws.Range("L7").Value = 3.95
$ws.Range("O7").Value = 1.39
$ws.Range("W7").Value = 6.3
$ws.Range("X7").Value = 9.75
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 22
$ws.Range("AA7").Value = 20
$ws.Range("AB7").Value = 35
$ws.Range("AG7").Value = 9
$ws.Range("AH7").Value = 18.5
$ws.Range("AI7").Value = 12
$ws.Range("AK7").Value = 35
$ws.Range("AL7").Value = 45
$ws.Range("AM7").Value = 700
$ws.Range("AN7").Value = 4
$ws.Range("AO7").Value = 11.5
$ws.Range("AP7").Value = 20
$ws.Range("AR7").Value = 80
$ws.Range("AX7").Value = 19.5
$ws.Range("AY7").Value = 25
$ws.Range("AZ7").Value = 100
$ws.Range("BA7").Value = 120
$ws.Range("BB7").Value = 300
$ws.Range("S8").Value = 1.5
$ws.Range("S9").Value = 1.5
$ws.Range("U9").Value = 1.87
$ws.Range("V9").Value = 1.77
$ws.Range("V10").Value = 1.58
$ws.Range("U11").Value = 1.69
$ws.Range("G12").Value = 1.27
$ws.Range("I12").Value = 12
$ws.Range("K12").Value = 2.4
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93
$ws.Range("U12").Value = 2.62
$ws.Range("V12").Value = 1.41
$ws.Range("Y12").Value = 10
$ws.Range("Z12").Value = 7
$ws.Range("AG12").Value = 21
$ws.Range("AI12").Value = 34
$ws.Range("AK12").Value = 101
$ws.Range("AL12").Value = 101
$ws.Range("AU12").Value = 12
$ws.Range("AZ12").Value = 351
